$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.390.73"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.639.12"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.531"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.77%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "1.871.97"
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("D13").Value = "1.646.34"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.555"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").Value = "27.356.13"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.20%  "
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  -4.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.77%  "
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.84%  "
$ws.Range("E31").Value = "  -3.86%  "
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "1.404.53"
$ws.Range("E34").Value = "  -4.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.560"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.11%  "
$ws.Range("E39").Value = "  -3.55%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.28%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.788"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "1.780.72"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("E48").Value = "  -4.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").Value = "0.0₆0104"
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0988"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.07%  "
